$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells whose new value is a plain decimal number must be forced to
# remain Text (matching the original inlineStr cell type) instead of being
# auto-converted to a Number by Excel's input parser. We set NumberFormat to
# "@" (Text) before assigning the value, then reset the style back to "Normal"
# afterwards so the cell keeps the workbook's original (unstyled) appearance.
$textForceRows = @(5, 6, 8, 10, 11, 16, 20, 21, 24, 28, 30, 33, 34, 35, 36, 40, 42, 43, 44, 45, 47, 48, 49, 50, 51)
foreach ($r in $textForceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "61.706.61"
$ws.Range("E2").Value = "  -4.44%  "

# Row 3
$ws.Range("D3").Value = "2.966.77"
$ws.Range("E3").Value = "  -6.40%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "540.28"
$ws.Range("E5").Value = "  -5.55%  "

# Row 6
$ws.Range("D6").Value = "152.74"
$ws.Range("E6").Value = "  -6.95%  "

# Row 7
$ws.Range("E7").Value = "  -0.17%  "

# Row 8
$ws.Range("D8").Value = "0.563"
$ws.Range("E8").Value = "  -4.05%  "

# Row 9
$ws.Range("D9").Value = "2.973.75"
$ws.Range("E9").Value = "  -6.06%  "

# Row 10
$ws.Range("D10").Value = "0.112"
$ws.Range("E10").Value = "  -4.52%  "

# Row 11
$ws.Range("D11").Value = "6.09"
$ws.Range("E11").Value = "  -8.22%  "

# Row 12
$ws.Range("E12").Value = "  -5.06%  "

# Row 13
$ws.Range("D13").Value = "3.482.04"
$ws.Range("E13").Value = "  -6.38%  "

# Row 14
$ws.Range("E14").Value = "  -3.64%  "

# Row 15
$ws.Range("D15").Value = "61.738.20"
$ws.Range("E15").Value = "  -4.40%  "

# Row 16
$ws.Range("D16").Value = "23.66"
$ws.Range("E16").Value = "  -6.60%  "

# Row 17
$ws.Range("D17").Value = "2.974.03"
$ws.Range("E17").Value = "  -6.27%  "

# Row 18
$ws.Range("E18").Value = "  -5.68%  "

# Row 19
$ws.Range("E19").Value = "  -2.90%  "

# Row 20
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "379.92"
$ws.Range("E20").Value = "  -7.13%  "

# Row 21
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "11.95"
$ws.Range("E21").Value = "  -6.28%  "

# Row 22
$ws.Range("E22").Value = "  -6.42%  "

# Row 23
$ws.Range("E23").Value = "  -0.02%  "

# Row 24
$ws.Range("D24").Value = "65.05"
$ws.Range("E24").Value = "  -5.62%  "

# Row 25
$ws.Range("E25").Value = "  -3.02%  "

# Row 26
$ws.Range("D26").Value = "3.093.23"
$ws.Range("E26").Value = "  -6.46%  "

# Row 27
$ws.Range("E27").Value = "  -5.16%  "

# Row 28
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.32%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0928"
$ws.Range("E29").Value = "  -9.41%  "

# Row 30
$ws.Range("D30").Value = "8.23"
$ws.Range("E30").Value = "  -7.27%  "

# Row 31
$ws.Range("E31").Value = "  -0.02%  "

# Row 32
$ws.Range("E32").Value = "  -5.83%  "

# Row 33
$ws.Range("D33").Value = "20.43"
$ws.Range("E33").Value = "  -3.93%  "

# Row 34
$ws.Range("D34").Value = "158.77"
$ws.Range("E34").Value = "  +1.55%  "

# Row 35
$ws.Range("D35").Value = "4.64"
$ws.Range("E35").Value = "  -5.32%  "

# Row 36
$ws.Range("D36").Value = "5.93"
$ws.Range("E36").Value = "  -6.66%  "

# Row 37
$ws.Range("E37").Value = "  -5.71%  "

# Row 38
$ws.Range("E38").Value = "  -5.65%  "

# Row 39
$ws.Range("E39").Value = "  -8.48%  "

# Row 40
$ws.Range("D40").Value = "3.91"
$ws.Range("E40").Value = "  -4.79%  "

# Row 41
$ws.Range("D41").Value = "2.413.95"
$ws.Range("E41").Value = "  -10.06%  "

# Row 42
$ws.Range("D42").Value = "37.10"
$ws.Range("E42").Value = "  -3.82%  "

# Row 43
$ws.Range("D43").Value = "22.22"
$ws.Range("E43").Value = "  -7.56%  "

# Row 44
$ws.Range("D44").Value = "0.663"
$ws.Range("E44").Value = "  -4.80%  "

# Row 45
$ws.Range("D45").Value = "0.0591"
$ws.Range("E45").Value = "  -4.74%  "

# Row 46
$ws.Range("E46").Value = "  -0.29%  "

# Row 47
$ws.Range("D47").Value = "0.0244"
$ws.Range("E47").Value = "  -5.49%  "

# Row 48
$ws.Range("D48").Value = "4.97"
$ws.Range("E48").Value = "  -8.68%  "

# Row 49
$ws.Range("D49").Value = "0.0948"
$ws.Range("E49").Value = "  -3.94%  "

# Row 50
$ws.Range("D50").Value = "266.44"
$ws.Range("E50").Value = "  -8.80%  "

# Row 51
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "10.47"
$ws.Range("E51").Value = "  +0.11%  "

# Reset style on the force-text cells back to Normal so no explicit style
# index lingers on the cell (keeps cells styleless, like the original file).
foreach ($r in $textForceRows) {
    $ws.Cells.Item($r, 4).Style = "Normal"
}
